$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the slightly adjusted timestamp value in A7
$ws.Range("A7").Value = 45878.25014258102

# Add new row 8 with data
$ws.Range("A8").Value = 45878.29184623543
$ws.Range("B8").Value = 2025
$ws.Range("C8").Value = 37
$ws.Range("D8").Value = 13
$ws.Range("E8").Value = 92.56
$ws.Range("F8").Value = 7.03
$ws.Range("G8").Value = 0.01
$ws.Range("H8").Value = "WNW"
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = "07:00:15"

# Apply same style as A2:A7 (style index 2, date/time format) to A8
$ws.Range("A8").NumberFormat = $ws.Range("A7").NumberFormat
